# edits to dickson protocol, wait 90 secs before starting experiment
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# temperature_setpoint: 22 -> 24
$ws.Range("B2").Value = 24

# gender value: "b" -> "m or f"
$ws.Range("B7").Value = "m or f"

# EP_version_n: 1 -> 1.1
$ws.Range("B12").Value = 1.1

# date_created: 20150330 -> 20180315
$ws.Range("B14").Value = 20180315

# New row 19, modeled on row 18's formatting (ht=22, customHeight, s=20/19/19)
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("D18").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 22
$ws.Range("A19").Value = "overhead lights off, screns dimmed with red channel lowered"

# selection moves to A20
$ws.Range("A20").Select() | Out-Null
